# Update "想去人数" (column F) values across the four worksheets to match
# the latest scrape snapshot referenced in the commit message.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    4  = 4352
    6  = 422
    7  = 3456
    8  = 981
    11 = 300
    12 = 2347
    16 = 504
    19 = 9737
    20 = 5926
    23 = 813
    25 = 826
    26 = 3512
    29 = 451
    30 = 104
    31 = 230
    32 = 213
    33 = 4779
    35 = 1036
    36 = 135
    37 = 16
    38 = 457
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(18, 6).Value = 5

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 8645
    4 = 1517
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 8645
    5  = 1517
    7  = 4352
    9  = 422
    10 = 3456
    11 = 981
    14 = 300
    15 = 2347
    23 = 504
    26 = 9737
    31 = 813
    33 = 826
    34 = 3512
    37 = 451
    38 = 104
    39 = 230
    40 = 5
    41 = 213
    42 = 4779
    43 = 1036
    44 = 135
    45 = 457
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
